# =====================================================================
# Update odds data for "Spain Segunda" worksheet
# Commit: Atualizacao de bases das ligas, do dia: 30-03-2024 as 19:32
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: refreshed odds data between paired rows (B:AC). ---
# --- Column A (sequential id) is untouched for each row number.   ---

# Rows 14 <-> 15
$ws.Range("B14").Value = 6837956
$ws.Range("B15").Value = 6839317
$ws.Range("F14").Value = "Albacete"
$ws.Range("F15").Value = "FC Cartagena"
$ws.Range("G14").Value = "Espanyol"
$ws.Range("G15").Value = "Eldense"
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("J14").Value = "D"
$ws.Range("J15").Value = "A"
$ws.Range("K14").Value = 2.875
$ws.Range("K15").Value = 1.909
$ws.Range("M14").Value = 2.25
$ws.Range("M15").Value = 3.75
$ws.Range("N14").Value = 2.25
$ws.Range("N15").Value = 2.2
$ws.Range("O14").Value = 3.3
$ws.Range("O15").Value = 3.1
$ws.Range("P14").Value = 3.3
$ws.Range("P15").Value = 3.6
$ws.Range("R14").Value = 1.9
$ws.Range("R15").Value = 1.85
$ws.Range("S14").Value = 1.95
$ws.Range("S15").Value = 2
$ws.Range("T14").Value = 2.25
$ws.Range("T15").Value = 2
$ws.Range("U14").Value = 1.95
$ws.Range("U15").Value = 1.9
$ws.Range("V14").Value = 1.9
$ws.Range("V15").Value = 1.95
$ws.Range("X14").Value = 2.3
$ws.Range("X15").Value = -1
$ws.Range("Y14").Value = -1
$ws.Range("Y15").Value = 2.6
$ws.Range("Z14").Value = -0.5
$ws.Range("Z15").Value = -1
$ws.Range("AA14").Value = 0.475
$ws.Range("AA15").Value = 1
$ws.Range("AB14").Value = -0.5
$ws.Range("AB15").Value = -1
$ws.Range("AC14").Value = 0.45
$ws.Range("AC15").Value = 0.95

# Rows 45 <-> 46
$ws.Range("B45").Value = 6838103
$ws.Range("B46").Value = 6838106
$ws.Range("F45").Value = "Elche"
$ws.Range("F46").Value = "Villarreal B"
$ws.Range("G45").Value = "Racing Santander"
$ws.Range("G46").Value = "FC Cartagena"
$ws.Range("I45").Value = 1
$ws.Range("I46").Value = 2
$ws.Range("J45").Value = "D"
$ws.Range("J46").Value = "A"
$ws.Range("K45").Value = 1.85
$ws.Range("K46").Value = 2.3
$ws.Range("L45").Value = 3.25
$ws.Range("L46").Value = 3.2
$ws.Range("M45").Value = 4.333
$ws.Range("M46").Value = 3.2
$ws.Range("N45").Value = 1.85
$ws.Range("N46").Value = 2.05
$ws.Range("P45").Value = 4.75
$ws.Range("P46").Value = 3.75
$ws.Range("R45").Value = 1.85
$ws.Range("R46").Value = 2.05
$ws.Range("S45").Value = 2
$ws.Range("S46").Value = 1.8
$ws.Range("T45").Value = 2.25
$ws.Range("T46").Value = 2.5
$ws.Range("U45").Value = 1.925
$ws.Range("U46").Value = 1.95
$ws.Range("V45").Value = 1.925
$ws.Range("V46").Value = 1.9
$ws.Range("X45").Value = 2.4
$ws.Range("X46").Value = -1
$ws.Range("Y45").Value = -1
$ws.Range("Y46").Value = 2.75
$ws.Range("AA45").Value = 1
$ws.Range("AA46").Value = 0.8
$ws.Range("AB45").Value = -0.5
$ws.Range("AB46").Value = 0.95
$ws.Range("AC45").Value = 0.4625
$ws.Range("AC46").Value = -1

# Rows 57 <-> 58
$ws.Range("B57").Value = 6839309
$ws.Range("B58").Value = 6836525
$ws.Range("F57").Value = "Racing Ferrol"
$ws.Range("F58").Value = "Racing Santander"
$ws.Range("G57").Value = "Villarreal B"
$ws.Range("G58").Value = "Amorebieta"
$ws.Range("H57").Value = 2
$ws.Range("H58").Value = 1
$ws.Range("I57").Value = 2
$ws.Range("I58").Value = 0
$ws.Range("J57").Value = "D"
$ws.Range("J58").Value = "H"
$ws.Range("K57").Value = 2.2
$ws.Range("K58").Value = 1.909
$ws.Range("L57").Value = 3
$ws.Range("L58").Value = 3.1
$ws.Range("M57").Value = 3.3
$ws.Range("M58").Value = 4.2
$ws.Range("N57").Value = 1.909
$ws.Range("N58").Value = 1.8
$ws.Range("O57").Value = 3.6
$ws.Range("O58").Value = 3.4
$ws.Range("P57").Value = 4
$ws.Range("P58").Value = 4.75
$ws.Range("Q57").Value = -0.5
$ws.Range("Q58").Value = -0.75
$ws.Range("R57").Value = 1.85
$ws.Range("R58").Value = 2.025
$ws.Range("S57").Value = 2
$ws.Range("S58").Value = 1.825
$ws.Range("U57").Value = 1.875
$ws.Range("U58").Value = 2
$ws.Range("V57").Value = 1.975
$ws.Range("V58").Value = 1.85
$ws.Range("W57").Value = -1
$ws.Range("W58").Value = 0.8
$ws.Range("X57").Value = 2.6
$ws.Range("X58").Value = -1
$ws.Range("Z57").Value = -1
$ws.Range("Z58").Value = 0.5125
$ws.Range("AA57").Value = 1
$ws.Range("AA58").Value = -0.5
$ws.Range("AB57").Value = 0.875
$ws.Range("AB58").Value = -1
$ws.Range("AC57").Value = -1
$ws.Range("AC58").Value = 0.8500000000000001

# Rows 120 <-> 121
$ws.Range("B120").Value = 6838121
$ws.Range("B121").Value = 6838006
$ws.Range("F120").Value = "FC Cartagena"
$ws.Range("F121").Value = "Sporting Gijon"
$ws.Range("G120").Value = "Racing Santander"
$ws.Range("G121").Value = "Real Zaragoza"
$ws.Range("I120").Value = 3
$ws.Range("I121").Value = 2
$ws.Range("J120").Value = "A"
$ws.Range("J121").Value = "D"
$ws.Range("K120").Value = 2.25
$ws.Range("K121").Value = 2.05
$ws.Range("M120").Value = 3.4
$ws.Range("M121").Value = 3.8
$ws.Range("N120").Value = 2.375
$ws.Range("N121").Value = 1.95
$ws.Range("O120").Value = 3
$ws.Range("O121").Value = 3.2
$ws.Range("P120").Value = 3.3
$ws.Range("P121").Value = 4.333
$ws.Range("Q120").Value = -0.25
$ws.Range("Q121").Value = -0.5
$ws.Range("R120").Value = 2
$ws.Range("R121").Value = 1.95
$ws.Range("S120").Value = 1.85
$ws.Range("S121").Value = 1.9
$ws.Range("U120").Value = 1.95
$ws.Range("U121").Value = 1.975
$ws.Range("V120").Value = 1.9
$ws.Range("V121").Value = 1.875
$ws.Range("X120").Value = -1
$ws.Range("X121").Value = 2.2
$ws.Range("Y120").Value = 2.3
$ws.Range("Y121").Value = -1
$ws.Range("AA120").Value = 0.8500000000000001
$ws.Range("AA121").Value = 0.8999999999999999
$ws.Range("AB120").Value = 0.95
$ws.Range("AB121").Value = 0.9750000000000001

# Rows 136 <-> 137
$ws.Range("B136").Value = 6838011
$ws.Range("B137").Value = 6838008
$ws.Range("F136").Value = "Racing Santander"
$ws.Range("F137").Value = "Albacete"
$ws.Range("G136").Value = "Burgos"
$ws.Range("G137").Value = "Sporting Gijon"
$ws.Range("H136").Value = 3
$ws.Range("H137").Value = 1
$ws.Range("I136").Value = 0
$ws.Range("I137").Value = 3
$ws.Range("J136").Value = "H"
$ws.Range("J137").Value = "A"
$ws.Range("M136").Value = 3.6
$ws.Range("M137").Value = 3.75
$ws.Range("N136").Value = 2.05
$ws.Range("N137").Value = 2.15
$ws.Range("O136").Value = 3.25
$ws.Range("O137").Value = 3.2
$ws.Range("P136").Value = 4
$ws.Range("P137").Value = 3.6
$ws.Range("Q136").Value = -0.5
$ws.Range("Q137").Value = -0.25
$ws.Range("R136").Value = 2.05
$ws.Range("R137").Value = 1.825
$ws.Range("S136").Value = 1.8
$ws.Range("S137").Value = 2.025
$ws.Range("U136").Value = 1.875
$ws.Range("U137").Value = 1.975
$ws.Range("V136").Value = 1.975
$ws.Range("V137").Value = 1.875
$ws.Range("W136").Value = 1.05
$ws.Range("W137").Value = -1
$ws.Range("Y136").Value = -1
$ws.Range("Y137").Value = 2.6
$ws.Range("Z136").Value = 1.05
$ws.Range("Z137").Value = -1
$ws.Range("AA136").Value = -1
$ws.Range("AA137").Value = 1.025
$ws.Range("AB136").Value = 0.875
$ws.Range("AB137").Value = 0.9750000000000001

# Rows 164 <-> 165
$ws.Range("B164").Value = 6839289
$ws.Range("B165").Value = 6838024
$ws.Range("F164").Value = "Mirandes"
$ws.Range("F165").Value = "Elche"
$ws.Range("G164").Value = "Racing Ferrol"
$ws.Range("G165").Value = "Real Zaragoza"
$ws.Range("H164").Value = 1
$ws.Range("H165").Value = 2
$ws.Range("I164").Value = 2
$ws.Range("I165").Value = 0
$ws.Range("J164").Value = "A"
$ws.Range("J165").Value = "H"
$ws.Range("K164").Value = 2.4
$ws.Range("K165").Value = 2.05
$ws.Range("L164").Value = 3.1
$ws.Range("L165").Value = 3.2
$ws.Range("M164").Value = 3.1
$ws.Range("M165").Value = 3.75
$ws.Range("N164").Value = 2.4
$ws.Range("N165").Value = 2.05
$ws.Range("O164").Value = 3.1
$ws.Range("O165").Value = 3.3
$ws.Range("P164").Value = 3.2
$ws.Range("P165").Value = 4
$ws.Range("Q164").Value = -0.25
$ws.Range("Q165").Value = -0.5
$ws.Range("R164").Value = 2
$ws.Range("R165").Value = 2.025
$ws.Range("S164").Value = 1.85
$ws.Range("S165").Value = 1.825
$ws.Range("T164").Value = 2
$ws.Range("T165").Value = 2.25
$ws.Range("U164").Value = 1.825
$ws.Range("U165").Value = 2.05
$ws.Range("V164").Value = 2.025
$ws.Range("V165").Value = 1.8
$ws.Range("W164").Value = -1
$ws.Range("W165").Value = 1.05
$ws.Range("Y164").Value = 2.2
$ws.Range("Y165").Value = -1
$ws.Range("Z164").Value = -1
$ws.Range("Z165").Value = 1.025
$ws.Range("AA164").Value = 0.8500000000000001
$ws.Range("AA165").Value = -1
$ws.Range("AB164").Value = 0.825
$ws.Range("AB165").Value = -0.5
$ws.Range("AC164").Value = -1
$ws.Range("AC165").Value = 0.4

# Rows 169 <-> 170
$ws.Range("B169").Value = 6838026
$ws.Range("B170").Value = 6838028
$ws.Range("F169").Value = "Tenerife"
$ws.Range("F170").Value = "Huesca"
$ws.Range("G169").Value = "Villarreal B"
$ws.Range("G170").Value = "Espanyol"
$ws.Range("H169").Value = 0
$ws.Range("H170").Value = 1
$ws.Range("J169").Value = "A"
$ws.Range("J170").Value = "D"
$ws.Range("K169").Value = 1.666
$ws.Range("K170").Value = 3.5
$ws.Range("L169").Value = 3.75
$ws.Range("L170").Value = 3.2
$ws.Range("M169").Value = 5
$ws.Range("M170").Value = 2.15
$ws.Range("N169").Value = 1.7
$ws.Range("N170").Value = 3.4
$ws.Range("O169").Value = 3.6
$ws.Range("O170").Value = 3.1
$ws.Range("P169").Value = 5.25
$ws.Range("P170").Value = 2.3
$ws.Range("Q169").Value = -0.75
$ws.Range("Q170").Value = 0.25
$ws.Range("R169").Value = 1.975
$ws.Range("R170").Value = 1.875
$ws.Range("S169").Value = 1.875
$ws.Range("S170").Value = 1.975
$ws.Range("T169").Value = 2.25
$ws.Range("T170").Value = 2
$ws.Range("U169").Value = 1.95
$ws.Range("U170").Value = 1.9
$ws.Range("V169").Value = 1.9
$ws.Range("V170").Value = 1.95
$ws.Range("X169").Value = -1
$ws.Range("X170").Value = 2.1
$ws.Range("Y169").Value = 4.25
$ws.Range("Y170").Value = -1
$ws.Range("Z169").Value = -1
$ws.Range("Z170").Value = 0.4375
$ws.Range("AA169").Value = 0.875
$ws.Range("AA170").Value = -0.5
$ws.Range("AB169").Value = -1
$ws.Range("AB170").Value = 0
$ws.Range("AC169").Value = 0.8999999999999999
$ws.Range("AC170").Value = -0

# Rows 228 <-> 229
$ws.Range("B228").Value = 6839278
$ws.Range("B229").Value = 6838147
$ws.Range("F228").Value = "Albacete"
$ws.Range("F229").Value = "Racing Santander"
$ws.Range("G228").Value = "Eldense"
$ws.Range("G229").Value = "FC Andorra"
$ws.Range("H228").Value = 1
$ws.Range("H229").Value = 2
$ws.Range("I228").Value = 1
$ws.Range("I229").Value = 0
$ws.Range("J228").Value = "D"
$ws.Range("J229").Value = "H"
$ws.Range("K228").Value = 1.615
$ws.Range("K229").Value = 2.35
$ws.Range("L228").Value = 3.75
$ws.Range("L229").Value = 3.2
$ws.Range("M228").Value = 5.5
$ws.Range("M229").Value = 3.1
$ws.Range("N228").Value = 1.65
$ws.Range("N229").Value = 2
$ws.Range("O228").Value = 3.8
$ws.Range("O229").Value = 3.5
$ws.Range("P228").Value = 5.5
$ws.Range("P229").Value = 3.75
$ws.Range("Q228").Value = -0.75
$ws.Range("Q229").Value = -0.5
$ws.Range("R228").Value = 1.825
$ws.Range("R229").Value = 2.025
$ws.Range("S228").Value = 2.025
$ws.Range("S229").Value = 1.825
$ws.Range("W228").Value = -1
$ws.Range("W229").Value = 1
$ws.Range("X228").Value = 2.8
$ws.Range("X229").Value = -1
$ws.Range("Z228").Value = -1
$ws.Range("Z229").Value = 1.025
$ws.Range("AA228").Value = 1.025
$ws.Range("AA229").Value = -1

# Rows 233 <-> 234
$ws.Range("B233").Value = 6838050
$ws.Range("B234").Value = 6838051
$ws.Range("F233").Value = "Eibar"
$ws.Range("F234").Value = "Leganes"
$ws.Range("G233").Value = "Sporting Gijon"
$ws.Range("G234").Value = "Tenerife"
$ws.Range("K233").Value = 1.85
$ws.Range("K234").Value = 2.05
$ws.Range("L233").Value = 3.4
$ws.Range("L234").Value = 3.2
$ws.Range("M233").Value = 4.2
$ws.Range("M234").Value = 4
$ws.Range("N233").Value = 1.909
$ws.Range("N234").Value = 2.3
$ws.Range("O233").Value = 3.5
$ws.Range("O234").Value = 2.8
$ws.Range("P233").Value = 4
$ws.Range("P234").Value = 3.6
$ws.Range("Q233").Value = -0.5
$ws.Range("Q234").Value = -0.25
$ws.Range("R233").Value = 1.95
$ws.Range("R234").Value = 1.975
$ws.Range("S233").Value = 1.9
$ws.Range("S234").Value = 1.875
$ws.Range("T233").Value = 2.25
$ws.Range("T234").Value = 1.75
$ws.Range("U233").Value = 1.975
$ws.Range("U234").Value = 2.05
$ws.Range("V233").Value = 1.875
$ws.Range("V234").Value = 1.8
$ws.Range("X233").Value = 2.5
$ws.Range("X234").Value = 1.8
$ws.Range("Z233").Value = -1
$ws.Range("Z234").Value = -0.5
$ws.Range("AA233").Value = 0.8999999999999999
$ws.Range("AA234").Value = 0.4375
$ws.Range("AB233").Value = -0.5
$ws.Range("AB234").Value = 0.5249999999999999
$ws.Range("AC233").Value = 0.4375
$ws.Range("AC234").Value = -0.5

# Rows 235 <-> 236
$ws.Range("B235").Value = 6838149
$ws.Range("B236").Value = 6838148
$ws.Range("F235").Value = "Villarreal B"
$ws.Range("F236").Value = "Huesca"
$ws.Range("G235").Value = "Oviedo"
$ws.Range("G236").Value = "FC Cartagena"
$ws.Range("H235").Value = 1
$ws.Range("H236").Value = 3
$ws.Range("I235").Value = 1
$ws.Range("I236").Value = 0
$ws.Range("J235").Value = "D"
$ws.Range("J236").Value = "H"
$ws.Range("K235").Value = 2.6
$ws.Range("K236").Value = 2.2
$ws.Range("L235").Value = 3.2
$ws.Range("L236").Value = 3.1
$ws.Range("M235").Value = 2.75
$ws.Range("M236").Value = 3.6
$ws.Range("N235").Value = 3
$ws.Range("N236").Value = 2.05
$ws.Range("O235").Value = 3
$ws.Range("O236").Value = 2.9
$ws.Range("P235").Value = 2.625
$ws.Range("P236").Value = 4.5
$ws.Range("Q235").Value = 0
$ws.Range("Q236").Value = -0.5
$ws.Range("T235").Value = 2
$ws.Range("T236").Value = 1.75
$ws.Range("W235").Value = -1
$ws.Range("W236").Value = 1.05
$ws.Range("X235").Value = 2
$ws.Range("X236").Value = -1
$ws.Range("Z235").Value = 0
$ws.Range("Z236").Value = 1.025
$ws.Range("AA235").Value = -0
$ws.Range("AA236").Value = -1
$ws.Range("AB235").Value = 0
$ws.Range("AB236").Value = 0.95
$ws.Range("AC235").Value = -0
$ws.Range("AC236").Value = -1

# --- Section 2: refresh odds data for upcoming (not-yet-played) fixtures ---
# Rows 360-363 get new id/date/teams and refreshed odds, matching data that
# used to live (with slightly different odds) in rows 366-369 respectively.

# Row 360
$ws.Range("B360").Value = 7129606
$ws.Range("E360").Value = 45382.375
$ws.Range("F360").Value = "Alcorcon"
$ws.Range("G360").Value = "Amorebieta"
$ws.Range("K360").Value = 2
$ws.Range("L360").Value = 3
$ws.Range("M360").Value = 4.333
$ws.Range("N360").Value = 2.15
$ws.Range("O360").Value = 3
$ws.Range("P360").Value = 4
$ws.Range("Q360").Value = -0.25
$ws.Range("R360").Value = 1.85
$ws.Range("S360").Value = 2
$ws.Range("T360").Value = 1.75
$ws.Range("U360").Value = 1.8
$ws.Range("V360").Value = 2.05

# Row 361
$ws.Range("B361").Value = 7129608
$ws.Range("E361").Value = 45382.46875
$ws.Range("F361").Value = "Real Zaragoza"
$ws.Range("G361").Value = "Tenerife"
$ws.Range("K361").Value = 2.3
$ws.Range("L361").Value = 2.75
$ws.Range("M361").Value = 3.75
$ws.Range("N361").Value = 2.45
$ws.Range("O361").Value = 2.6
$ws.Range("P361").Value = 3.75
$ws.Range("Q361").Value = -0.25
$ws.Range("R361").Value = 2
$ws.Range("S361").Value = 1.85
$ws.Range("T361").Value = 1.75
$ws.Range("U361").Value = 2.05
$ws.Range("V361").Value = 1.8

# Row 362
$ws.Range("B362").Value = 7129607
$ws.Range("E362").Value = 45382.46875
$ws.Range("F362").Value = "Racing Ferrol"
$ws.Range("G362").Value = "Elche"
$ws.Range("K362").Value = 2.5
$ws.Range("L362").Value = 2.8
$ws.Range("M362").Value = 3.2
$ws.Range("N362").Value = 3.1
$ws.Range("O362").Value = 2.9
$ws.Range("P362").Value = 2.625
$ws.Range("Q362").Value = 0
$ws.Range("R362").Value = 2.1
$ws.Range("S362").Value = 1.775
$ws.Range("T362").Value = 2
$ws.Range("U362").Value = 1.95
$ws.Range("V362").Value = 1.9

# Row 363
$ws.Range("B363").Value = 7128454
$ws.Range("E363").Value = 45382.5625
$ws.Range("F363").Value = "Oviedo"
$ws.Range("G363").Value = "Villarreal B"
$ws.Range("K363").Value = 1.533
$ws.Range("L363").Value = 4
$ws.Range("M363").Value = 6
$ws.Range("N363").Value = 1.533
$ws.Range("O363").Value = 4
$ws.Range("P363").Value = 6.5
$ws.Range("Q363").Value = -1
$ws.Range("R363").Value = 1.925
$ws.Range("S363").Value = 1.925
$ws.Range("T363").Value = 2.5
$ws.Range("U363").Value = 1.975
$ws.Range("V363").Value = 1.875

# --- Section 3: remove the now-obsolete trailing fixtures (former rows 364-369) ---
$ws.Range("A364:AC369").EntireRow.Delete() | Out-Null
